$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look like numbers are written as plain text,
# matching the source data which stores localized price strings (e.g. "1.007", "27.087.84").

$ws.Range('D2').Value = '27.087.84'
$ws.Range('E2').Value = '  -1.56%  '

$ws.Range('D3').Value = '1.798.21'
$ws.Range('E3').Value = '  -2.31%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.40%  '

$ws.Range('E5').Value = '  +0.34%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '307.81'
$ws.Range('E6').Value = '  -1.94%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4175'
$ws.Range('E7').Value = '  -1.90%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3558'
$ws.Range('E8').Value = '  -3.04%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07062'
$ws.Range('E9').Value = '  -3.12%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8434'
$ws.Range('E10').Value = '  -3.20%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.23'
$ws.Range('E11').Value = '  -2.89%  '

$ws.Range('D12').Value = '1.741.14'
$ws.Range('E12').Value = '  -6.86%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.267'
$ws.Range('E13').Value = '  -2.63%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.338'
$ws.Range('E14').Value = '  -2.98%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.06789'

$ws.Range('E16').Value = '  +0.52%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '79.92'
$ws.Range('E17').Value = '  -0.59%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008702'
$ws.Range('E18').Value = '  -3.54%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.007'
$ws.Range('E19').Value = '  +0.44%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.07'
$ws.Range('E20').Value = '  -2.69%  '

$ws.Range('D21').Value = '27.362.91'
$ws.Range('E21').Value = '  -1.12%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.051'
$ws.Range('E22').Value = '  -0.43%  '

$ws.Range('E23').Value = '  -0.80%  '

$ws.Range('D24').Value = '2.077.12'
$ws.Range('E24').Value = '  -0.22%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.955'
$ws.Range('E25').Value = '  -0.28%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '153.28'
$ws.Range('E26').Value = '  -0.82%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.10'
$ws.Range('E27').Value = '  -1.76%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.021'
$ws.Range('E28').Value = '  -4.46%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '112.72'
$ws.Range('E29').Value = '  -2.36%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.656'
$ws.Range('E30').Value = '  -11.44%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08876'
$ws.Range('E31').Value = '  -0.15%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7190'
$ws.Range('E32').Value = '  -7.79%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.867'
$ws.Range('E33').Value = '  -3.35%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.339'
$ws.Range('E34').Value = '  -4.98%  '

$ws.Range('E35').Value = '  +0.38%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.074'
$ws.Range('E36').Value = '  -6.96%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.076'
$ws.Range('E37').Value = '  -2.42%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01892'
$ws.Range('E38').Value = '  -3.08%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05101'
$ws.Range('E39').Value = '  -5.50%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.1617'
$ws.Range('E40').Value = '  -2.90%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4935'
$ws.Range('E41').Value = '  -4.01%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.594'
$ws.Range('E42').Value = '  -8.05%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.168'
$ws.Range('E43').Value = '  -8.83%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.045'
$ws.Range('E44').Value = '  -6.37%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.006'
$ws.Range('E45').Value = '  +0.34%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '104.45'
$ws.Range('E46').Value = '  -2.15%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.19'
$ws.Range('E47').Value = '  -3.14%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.06311'
$ws.Range('E48').Value = '  -3.41%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.4509'
$ws.Range('E49').Value = '  -4.84%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.585'
$ws.Range('E50').Value = '  -3.53%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '61.96'
$ws.Range('E51').Value = '  -4.01%  '

